# Auto-save via app Streamlit
# A new reservation row (Claudine Fleury / Booking) is inserted right above the
# existing "TOTAL" row, pushing that TOTAL row from row 51 down to row 52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 51 -- this shifts the old row 51 (TOTAL)
# down to row 52, and any formatting on row 51 (date styles in D/E) moves
# along with it automatically.
$ws.Rows.Item(51).Insert()

# Fill in the new row 51 with the new reservation's data (same values as the
# two existing "Claudine Fleury" rows directly above it, rows 49-50).
$ws.Cells.Item(51, 1).Value = "Claudine Fleury"
$ws.Cells.Item(51, 2).Value = "Booking"

# Phone number is stored as literal text (with leading "+"), not a number.
$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = "+33334691787"

# Arrival / departure dates (serial date values), formatted as dates like the
# cells above them.
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(51, 4).Value = 46200
$ws.Cells.Item(51, 5).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(51, 5).Value = 46203

$ws.Cells.Item(51, 6).Value = 3
$ws.Cells.Item(51, 7).Value = 468.68
$ws.Cells.Item(51, 8).Value = 374.77
$ws.Cells.Item(51, 9).Value = 93.91
$ws.Cells.Item(51, 10).Value = 20.04
$ws.Cells.Item(51, 11).Value = 2026
$ws.Cells.Item(51, 12).Value = 6
